# [ISP-2] Nuovi attributi Profondità e Fourilista
# Add a new row (16) to the Library_Formula sheet, replicating the
# "PuntoDiPrelievo" TS entries but for the new "TSLength" formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Library_Formula")

# Populate the new row's values (mirrors rows 2-7 / 14-15 pattern)
$ws.Range("A16").Value = "CREATE/MODIFY"
$ws.Range("B16").Value = "LibDemandPlanning"
$ws.Range("C16").Value = "TSLength"
$ws.Range("E16").Value = "TS"
$ws.Range("F16").Value = "PuntoDiPrelievo"

# Match the formatting used by the rest of the table: columns A, B and E
# carry the "body" font (same as row 15), while C and F keep the sheet's
# default column formatting.
foreach ($col in @("A", "B", "E")) {
    $src = $ws.Range($col + "15")
    $dst = $ws.Range($col + "16")
    $dst.Font.Name = $src.Font.Name
    $dst.Font.Size = $src.Font.Size
    $dst.Font.Color = $src.Font.Color
}

# Update the sheet's active selection to track the newly added row, as
# happens naturally when a user edits cell F16 last.
[void]$ws.Range("E15:F16").Select()
